# "include feels as hnames" -- append recode-book rows that map the
# Q3.32 / Q3.33 "feelings" battery sub-questions to their human-readable
# names (and alt-qnames), onto the end of the Sheet1 table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# qname, hname, altqname triples -> rows 21..32
$rows = @(
    @("Q3.33_1",  "Anxious",     "Q3.32_1"),
    @("Q3.33_2",  "Excited",     "Q3.32_2"),
    @("Q3.33_3",  "Frustrated",  "Q3.32_3"),
    @("Q3.33_4",  "Happy",       "Q3.32_4"),
    @("Q3.33_5",  "Scared",      "Q3.32_5"),
    @("Q3.33_11", "Ambivalent",  "Q3.32_11"),
    @("Q3.33_6",  "Surprised",   "Q3.32_6"),
    @("Q3.33_7",  "Thankful",    "Q3.32_7"),
    @("Q3.33_8",  "Unhappy",     "Q3.32_8"),
    @("Q3.33_9",  "Worried",     "Q3.32_9"),
    @("Q3.33_17", "Angry",       "Q3.32_17"),
    @("Q3.33_10", "Other",       "Q3.32_10")
)

$r = 21
foreach ($row in $rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Column widths: column A gets the (implicit) default width pinned
# explicitly, column B narrows slightly to fit the new content.
$ws1.Columns.Item(1).ColumnWidth = 7.671768707482998
$ws1.Columns.Item(2).ColumnWidth = 12.396258503401366

# Matches the recorded selection/active cell after the edit.
$ws1.Range("D23").Select()

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Columns.Item(1).ColumnWidth = 7.671768707482998
$ws2.Columns.Item(2).ColumnWidth = 27.110544217687067
